# Aggiornamento fino a 02/05: append new daily rows (239-244) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
  @(44313, 0, 11, 133.6898395721925),
  @(44314, 2, 12, 145.8434613514827),
  @(44315, 1, 13, 157.997083130773),
  @(44316, 1, 13, 157.997083130773),
  @(44317, 1, 14, 170.1507049100632),
  @(44318, 1, 12, 145.8434613514827)
)

$startRow = 239
for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $startRow + $i
  $row = $data[$i]

  # Copy the formatting (date style, border, bold, alignment) from the last
  # existing data row (238) onto the new row's date cell (column A) so the
  # style matches the rest of the column without introducing new style defs.
  $ws.Range("A238").Copy()
  $ws.Range("A$r").PasteSpecial(-4122)

  $ws.Range("A$r").Value = $row[0]
  $ws.Range("B$r").Value = $row[1]
  $ws.Range("C$r").Value = $row[2]
  $ws.Range("D$r").Value = $row[3]
}

$excel.CutCopyMode = 0
